# mosip_master/xlsx/role_list.xlsx - add REGISTRATION_OPERATOR role row
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "eng"
$ws.Range("B6").Value = "REGISTRATION_OPERATOR"
$ws.Range("C6").Value = "Registration Operator"

# Copy D2 (an existing "TRUE" text cell with the correct text-number-format
# style) into D6 so the new is_active value is written as the shared text
# string "TRUE" using the same style, rather than Excel's automatic boolean
# coercion.
$ws.Range("D2").Copy()
$ws.Range("D6").PasteSpecial()
$excel.CutCopyMode = $false

$ws.Range("D6").Select()
